$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 121; this shifts existing rows 121:137 down to 122:138
$ws.Rows("121:121").Insert()

# Populate the new row 121 with data (copy of former row 121 with updated price/date fields)
$ws.Range("A121").Value = 1
$ws.Range("B121").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C121").Value = "Arica y Parinacota"
$ws.Range("D121").Value = 44474
$ws.Range("E121").Value = 15
$ws.Range("F121").Value = "Fruta"
$ws.Range("G121").Value = 100108
$ws.Range("H121").Value = "Tropicales y subtropicales"
$ws.Range("I121").Value = 100108006
$ws.Range("J121").Value = "Plátano"
$ws.Range("K121").Value = "Sin especificar"
$ws.Range("L121").Value = "Pintón"
$ws.Range("M121").Value = 120
$ws.Range("N121").Value = 21000
$ws.Range("O121").Value = 22000
$ws.Range("P121").Value = 21500
$ws.Range("Q121").Value = "$/caja 20 kilos"
$ws.Range("R121").Value = "Ecuador"
$ws.Range("S121").Value = 1075
$ws.Range("T121").Value = 20
